$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 292.41177
$ws.Range("I33").Value = 205.2
$ws.Range("J33").Value = 417
$ws.Range("K33").Value = 205.2
$ws.Range("L33").Value = 417
$ws.Range("M33").Value = 23.80000000000001
$ws.Range("N33").Value = -875

$ws.Range("H132").Value = 1693290.8
$ws.Range("I132").Value = 3446
$ws.Range("K132").Value = 10338
$ws.Range("M132").Value = -7808

$ws.Range("H135").Value = 37030.035
$ws.Range("I135").Value = 42107.4
$ws.Range("J135").Value = 5296.5
$ws.Range("K135").Value = 378966.6
$ws.Range("L135").Value = 47668.5
$ws.Range("M135").Value = -376431.6
$ws.Range("N135").Value = -52738.5

$ws.Range("H137").Value = 5004105.5
$ws.Range("I137").Value = 8335996.5
$ws.Range("J137").Value = 6269.125
$ws.Range("K137").Value = 25007989.5
$ws.Range("L137").Value = 18807.375
$ws.Range("M137").Value = -25005439.5
$ws.Range("N137").Value = -23907.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 38540308
$ws.Range("I61").Value = 47667988
$ws.Range("J61").Value = 204060
$ws.Range("K61").Value = 47667988
$ws.Range("L61").Value = 204060
$ws.Range("M61").Value = -47667776
$ws.Range("N61").Value = -204484

$ws.Range("H74").Value = 6450593
$ws.Range("I74").Value = 11145802
$ws.Range("K74").Value = 11145802
$ws.Range("M74").Value = -11144928

$ws.Range("H77").Value = 6450593
$ws.Range("I77").Value = 11145802
$ws.Range("K77").Value = 55729010
$ws.Range("M77").Value = -55724642

$ws.Range("H132").Value = 58335.89
$ws.Range("I132").Value = 34206.934
$ws.Range("J132").Value = 183002.17
$ws.Range("K132").Value = 102620.802
$ws.Range("L132").Value = 549006.51
$ws.Range("M132").Value = -100090.802
$ws.Range("N132").Value = -554066.51

$ws.Range("H136").Value = 38540308
$ws.Range("I136").Value = 47667988
$ws.Range("J136").Value = 204060
$ws.Range("K136").Value = 143003964
$ws.Range("L136").Value = 612180
$ws.Range("M136").Value = -143001414
$ws.Range("N136").Value = -617280

$ws.Range("H138").Value = 52700
$ws.Range("J138").Value = 52700
$ws.Range("L138").Value = 52700
$ws.Range("N138").Value = -62980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3536.2856
$ws.Range("I134").Value = 3536.2856
$ws.Range("K134").Value = 10608.8568
$ws.Range("M134").Value = -8073.856800000001

$ws.Range("H140").Value = 48190
$ws.Range("J140").Value = 48190
$ws.Range("L140").Value = 48190
$ws.Range("N140").Value = -58550

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3431.0322
$ws.Range("I31").Value = 3940.182
$ws.Range("J31").Value = 3151
$ws.Range("K31").Value = 3940.182
$ws.Range("L31").Value = 3151
$ws.Range("M31").Value = -3645.182
$ws.Range("N31").Value = -3741

$ws.Range("H34").Value = 3431.0322
$ws.Range("I34").Value = 3940.182
$ws.Range("J34").Value = 3151
$ws.Range("K34").Value = 3940.182
$ws.Range("L34").Value = 3151
$ws.Range("M34").Value = -3738.182
$ws.Range("N34").Value = -3555

$ws.Range("H58").Value = 34484440
$ws.Range("I58").Value = 37038696
$ws.Range("K58").Value = 37038696
$ws.Range("M58").Value = -37038493

$ws.Range("H98").Value = 47138
$ws.Range("J98").Value = 47138
$ws.Range("L98").Value = 47138
$ws.Range("N98").Value = -51630

$ws.Range("H99").Value = 1275.28
$ws.Range("I99").Value = 1176.7778
$ws.Range("K99").Value = 1176.7778
$ws.Range("M99").Value = 321.2221999999999

$ws.Range("H100").Value = 39990
$ws.Range("J100").Value = 39990
$ws.Range("L100").Value = 39990
$ws.Range("N100").Value = -42154

$ws.Range("H126").Value = 1275.28
$ws.Range("I126").Value = 1176.7778
$ws.Range("K126").Value = 3530.3334
$ws.Range("M126").Value = -1060.3334

$ws.Range("H127").Value = 31071.818
$ws.Range("J127").Value = 31071.818
$ws.Range("L127").Value = 31071.818
$ws.Range("N127").Value = -40991.818

$ws.Range("H133").Value = 26474.076
$ws.Range("I133").Value = 38000
$ws.Range("J133").Value = 26013.04
$ws.Range("K133").Value = 38000
$ws.Range("L133").Value = 26013.04
$ws.Range("M133").Value = -35470
$ws.Range("N133").Value = -31073.04

$ws.Range("H136").Value = 34484440
$ws.Range("I136").Value = 37038696
$ws.Range("K136").Value = 111116088
$ws.Range("M136").Value = -111113538

$ws.Range("H140").Value = 50780
$ws.Range("J140").Value = 50780
$ws.Range("L140").Value = 50780
$ws.Range("N140").Value = -61140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 3167.5
$ws.Range("I119").Value = 1316.1111
$ws.Range("J119").Value = 6500
$ws.Range("K119").Value = 3948.3333
$ws.Range("L119").Value = 19500
$ws.Range("M119").Value = 889.6666999999998
$ws.Range("N119").Value = -29176

$ws.Range("H131").Value = 932.5714
$ws.Range("J131").Value = 1015.26086
$ws.Range("L131").Value = 3045.78258
$ws.Range("N131").Value = -13125.78258

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1697.6
$ws.Range("I122").Value = 1473.7778
$ws.Range("J122").Value = 2033.3334
$ws.Range("K122").Value = 4421.3334
$ws.Range("L122").Value = 6100.0002
$ws.Range("M122").Value = -1971.3334
$ws.Range("N122").Value = -11000.0002

$ws.Range("H132").Value = 169680
$ws.Range("I132").Value = 114018
$ws.Range("K132").Value = 342054
$ws.Range("M132").Value = -339524

$ws.Range("H133").Value = 47676.125
$ws.Range("J133").Value = 48671.43
$ws.Range("L133").Value = 48671.43
$ws.Range("N133").Value = -58791.43

$ws.Range("H135").Value = 36211.58
$ws.Range("J135").Value = 36211.58
$ws.Range("L135").Value = 36211.58
$ws.Range("N135").Value = -46351.58

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1440.8518
$ws.Range("I100").Value = 1123.7059
$ws.Range("J100").Value = 1980
$ws.Range("K100").Value = 1123.7059
$ws.Range("L100").Value = 1980
$ws.Range("M100").Value = -582.7058999999999
$ws.Range("N100").Value = -3062

$ws.Range("H132").Value = 50727.906
$ws.Range("I132").Value = 2482.3333
$ws.Range("J132").Value = 115055.336
$ws.Range("K132").Value = 7446.999899999999
$ws.Range("L132").Value = 345166.008
$ws.Range("M132").Value = -4916.999899999999
$ws.Range("N132").Value = -350226.008

$ws.Range("H134").Value = 59699.5
$ws.Range("J134").Value = 59699.5
$ws.Range("L134").Value = 59699.5
$ws.Range("N134").Value = -69839.5

$ws.Range("H136").Value = 224767.67
$ws.Range("I136").Value = 201400.8
$ws.Range("J136").Value = 253976.25
$ws.Range("K136").Value = 604202.3999999999
$ws.Range("L136").Value = 761928.75
$ws.Range("M136").Value = -601652.3999999999
$ws.Range("N136").Value = -767028.75

$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("M141").Value = -60360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8092.3335
$ws.Range("I41").Value = 9900
$ws.Range("J41").Value = 7188.5
$ws.Range("K41").Value = 9900
$ws.Range("L41").Value = 7188.5
$ws.Range("M41").Value = -9510
$ws.Range("N41").Value = -7968.5

$ws.Range("H45").Value = 8847.223
$ws.Range("J45").Value = 8953.125
$ws.Range("L45").Value = 8953.125
$ws.Range("N45").Value = -9935.125

$ws.Range("H74").Value = 5278.75
$ws.Range("J74").Value = 5278.75
$ws.Range("L74").Value = 5278.75
$ws.Range("N74").Value = -7150.75

$ws.Range("H77").Value = 5278.75
$ws.Range("J77").Value = 5278.75
$ws.Range("L77").Value = 15836.25
$ws.Range("N77").Value = -25196.25
